$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 66
$ws.Range("I31").Value = 66
$ws.Range("K31").Value = 198
$ws.Range("M31").Value = 32
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H51").Value = 7998.5
$ws.Range("J51").Value = 7998.5
$ws.Range("L51").Value = 7998.5
$ws.Range("N51").Value = -8966.5
$ws.Range("H52").Value = 3999.6667
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H82").Value = 394.33334
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 394.33334
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H118").Value = 363.2
$ws.Range("I118").Value = 376.75
$ws.Range("J118").Value = 309
$ws.Range("K118").Value = 1130.25
$ws.Range("L118").Value = 927
$ws.Range("M118").Value = 526.75
$ws.Range("N118").Value = -4241
$ws.Range("H129").Value = 15677.75
$ws.Range("I129").Value = 14164.667
$ws.Range("K129").Value = 42494.001
$ws.Range("M129").Value = -37494.001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 107.375
$ws.Range("I5").Value = 51.4
$ws.Range("K5").Value = 51.4
$ws.Range("M5").Value = 60.6
$ws.Range("H32").Value = 1994.5186
$ws.Range("I32").Value = 1994.5186
$ws.Range("K32").Value = 1994.5186
$ws.Range("M32").Value = -1707.5186
$ws.Range("H61").Value = 3186
$ws.Range("I61").Value = 2375
$ws.Range("J61").Value = 3997
$ws.Range("K61").Value = 2375
$ws.Range("L61").Value = 3997
$ws.Range("M61").Value = -2163
$ws.Range("N61").Value = -4421
$ws.Range("H74").Value = 1950
$ws.Range("I74").Value = 1900
$ws.Range("K74").Value = 1900
$ws.Range("M74").Value = -1026
$ws.Range("H77").Value = 1950
$ws.Range("I77").Value = 1900
$ws.Range("K77").Value = 9500
$ws.Range("M77").Value = -5132
$ws.Range("H136").Value = 3186
$ws.Range("I136").Value = 2375
$ws.Range("J136").Value = 3997
$ws.Range("K136").Value = 7125
$ws.Range("L136").Value = 11991
$ws.Range("M136").Value = -4575
$ws.Range("N136").Value = -17091
$ws.Range("H138").Value = 49999.5
$ws.Range("J138").Value = 49999.5
$ws.Range("L138").Value = 49999.5
$ws.Range("N138").Value = -60279.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 107.375
$ws.Range("I4").Value = 51.4
$ws.Range("K4").Value = 51.4
$ws.Range("M4").Value = 63.6
$ws.Range("H94").Value = 3307.2727
$ws.Range("I94").Value = 1876
$ws.Range("J94").Value = 4500
$ws.Range("K94").Value = 1876
$ws.Range("L94").Value = 4500
$ws.Range("M94").Value = -1425
$ws.Range("N94").Value = -5402

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1826.25
$ws.Range("I7").Value = 1438.625
$ws.Range("K7").Value = 1438.625
$ws.Range("M7").Value = -1325.625
$ws.Range("H31").Value = 2545.7778
$ws.Range("I31").Value = 1864
$ws.Range("K31").Value = 1864
$ws.Range("M31").Value = -1569
$ws.Range("H34").Value = 2545.7778
$ws.Range("I34").Value = 1864
$ws.Range("K34").Value = 1864
$ws.Range("M34").Value = -1662
$ws.Range("H58").Value = 2775.889
$ws.Range("I58").Value = 2333
$ws.Range("J58").Value = 2997.3333
$ws.Range("K58").Value = 2333
$ws.Range("L58").Value = 2997.3333
$ws.Range("M58").Value = -2130
$ws.Range("N58").Value = -3403.3333
$ws.Range("H99").Value = 8261.691999999999
$ws.Range("I99").Value = 6580.4
$ws.Range("J99").Value = 9312.5
$ws.Range("K99").Value = 6580.4
$ws.Range("L99").Value = 9312.5
$ws.Range("M99").Value = -5082.4
$ws.Range("N99").Value = -12308.5
$ws.Range("H105").Value = 2919.2
$ws.Range("I105").Value = 1774
$ws.Range("K105").Value = 1774
$ws.Range("M105").Value = -27
$ws.Range("H126").Value = 8261.691999999999
$ws.Range("I126").Value = 6580.4
$ws.Range("J126").Value = 9312.5
$ws.Range("K126").Value = 19741.2
$ws.Range("L126").Value = 27937.5
$ws.Range("M126").Value = -17271.2
$ws.Range("N126").Value = -32877.5
$ws.Range("H136").Value = 2775.889
$ws.Range("I136").Value = 2333
$ws.Range("J136").Value = 2997.3333
$ws.Range("K136").Value = 6999
$ws.Range("L136").Value = 8991.999899999999
$ws.Range("M136").Value = -4449
$ws.Range("N136").Value = -14091.9999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 35
$ws.Range("I26").Value = 35
$ws.Range("K26").Value = 105
$ws.Range("M26").Value = 183
$ws.Range("H69").Value = 2400
$ws.Range("I69").Value = 2400
$ws.Range("K69").Value = 7200
$ws.Range("M69").Value = -6389
$ws.Range("H72").Value = 2400
$ws.Range("I72").Value = 2400
$ws.Range("K72").Value = 21600
$ws.Range("M72").Value = -17544

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("H94").Value = 57000
$ws.Range("J94").Value = 57000
$ws.Range("L94").Value = 57000
$ws.Range("N94").Value = -58352
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H132").Value = 3082.2917
$ws.Range("I132").Value = 2999.3914
$ws.Range("K132").Value = 8998.174199999999
$ws.Range("M132").Value = -6468.174199999999
$ws.Range("H139").Value = 59999.332
$ws.Range("J139").Value = 59999.332
$ws.Range("L139").Value = 59999.332
$ws.Range("N139").Value = -70279.33199999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 662.13336
$ws.Range("I55").Value = 60.75
$ws.Range("J55").Value = 880.8182
$ws.Range("K55").Value = 60.75
$ws.Range("L55").Value = 880.8182
$ws.Range("M55").Value = 112.25
$ws.Range("N55").Value = -1226.8182
$ws.Range("H61").Value = 2290.5
$ws.Range("I61").Value = 2082
$ws.Range("J61").Value = 3750
$ws.Range("K61").Value = 2082
$ws.Range("L61").Value = 3750
$ws.Range("M61").Value = -1880
$ws.Range("N61").Value = -4154
$ws.Range("H113").Value = 2290.5
$ws.Range("I113").Value = 2082
$ws.Range("J113").Value = 3750
$ws.Range("K113").Value = 2082
$ws.Range("L113").Value = 3750
$ws.Range("M113").Value = 88
$ws.Range("N113").Value = -8090

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 42999
$ws.Range("J50").Value = 42999
$ws.Range("L50").Value = 42999
$ws.Range("N50").Value = -44261
$ws.Range("H119").Value = 72447.75
$ws.Range("J119").Value = 72447.75
$ws.Range("L119").Value = 72447.75
$ws.Range("N119").Value = -82123.75
$ws.Range("H138").Value = 40428
$ws.Range("J138").Value = 40428
$ws.Range("L138").Value = 40428
$ws.Range("N138").Value = -50708
